# Auto-generated Excel COM-interop script to update betting odds values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 2.88
$ws.Range("K2").Value = 1.95
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 2.4
$ws.Range("T2").Value = 1.53
$ws.Range("U2").Value = 3.95
$ws.Range("V2").Value = 1.26
$ws.Range("W2").Value = 4.5
$ws.Range("X2").Value = 1.18
$ws.Range("Y2").Value = 1.53
$ws.Range("Z2").Value = 2.38
$ws.Range("AH2").Value = 34
$ws.Range("AI2").Value = 7
$ws.Range("AK2").Value = 19
$ws.Range("AP2").Value = 13
# Row 3
$ws.Range("G3").Value = 2.25
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 1.8
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("Z3").Value = 2
$ws.Range("AD3").Value = 8.5
$ws.Range("AF3").Value = 21
$ws.Range("AG3").Value = 26
$ws.Range("AO3").Value = 17
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 41
$ws.Range("AS3").Value = 51
# Row 4
$ws.Range("G4").Value = 2.45
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 3.3
$ws.Range("K4").Value = 1.91
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 2.03
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 2.7
$ws.Range("T4").Value = 1.44
$ws.Range("U4").Value = 4.4
$ws.Range("Y4").Value = 1.62
$ws.Range("Z4").Value = 2.2
$ws.Range("AA4").Value = 2.1
$ws.Range("AB4").Value = 1.67
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 10
$ws.Range("AK4").Value = 19
$ws.Range("AP4").Value = 13
$ws.Range("AR4").Value = 34
# Row 5
$ws.Range("G5").Value = 2.88
$ws.Range("K5").Value = 1.73
$ws.Range("U5").Value = 7.6
$ws.Range("V5").Value = 1.09
$ws.Range("W5").Value = 10
$ws.Range("X5").Value = 1.06
$ws.Range("Y5").Value = 1.9
$ws.Range("Z5").Value = 1.9
$ws.Range("AD5").Value = 11
$ws.Range("AI5").Value = 4.5
$ws.Range("AK5").Value = 29
$ws.Range("AL5").Value = 151
# Row 7
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
# Row 8
$ws.Range("G8").Value = 2.5
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3.4
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("AC8").Value = 8
$ws.Range("AD8").Value = 12
$ws.Range("AF8").Value = 26
$ws.Range("AH8").Value = 34
$ws.Range("AO8").Value = 13
# Row 9
$ws.Range("S9").Value = 2.35
$ws.Range("T9").Value = 1.57
$ws.Range("W9").Value = 4.33
$ws.Range("X9").Value = 1.2
# Row 13
$ws.Range("G13").Value = 1.95
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 3.7
$ws.Range("J13").Value = 2.75
$ws.Range("K13").Value = 1.91
$ws.Range("L13").Value = 5
$ws.Range("Q13").Value = 2.03
$ws.Range("R13").Value = 1.83
$ws.Range("Y13").Value = 1.62
$ws.Range("Z13").Value = 2.2
$ws.Range("AA13").Value = 2.38
$ws.Range("AB13").Value = 1.53
$ws.Range("AD13").Value = 8
$ws.Range("AF13").Value = 17
$ws.Range("AK13").Value = 23
$ws.Range("AN13").Value = 8
$ws.Range("AO13").Value = 19
$ws.Range("AP13").Value = 15
$ws.Range("AR13").Value = 41
$ws.Range("AS13").Value = 51
# Row 14
$ws.Range("T14").Value = 1.67
